# Generate Report for Handback
# This script fills in the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns on the per-language sheets (zh-cn, de-de)
# now that a handback has completed, and updates the Overview sheet's status
# text accordingly.

$wb = $excel.ActiveWorkbook

$urlMd1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bf4b80184046b22466cdebd8d46a68a481d41ff5/e2e/07f5c7a9-97d9-4a11-9d2a-0d22f6f50138.md"
$urlMd2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bf4b80184046b22466cdebd8d46a68a481d41ff5/e2e/7f7f1aed-0a42-45ff-a0aa-1e0ff486b802.md"

$name1 = "07f5c7a9-97d9-4a11-9d2a-0d22f6f50138.md"
$name2 = "7f7f1aed-0a42-45ff-a0aa-1e0ff486b802.md"

# ---------------------------------------------------------------------
# Overview sheet: status columns (zh-cn / de-de) now report handback done
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("I2").Value = $name1
$wsZh.Range("J2").Value = "07f5c7a9-97d9-4a11-9d2a-0d22f6f50138.bdc57f7efbf1e5bb7004563fe238e320ad0e99a9.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-09-03 23:05:46"

$wsZh.Range("I3").Value = $name2
$wsZh.Range("J3").Value = "7f7f1aed-0a42-45ff-a0aa-1e0ff486b802.6f682d1f368a32d8bd4f8857abd3b3d30eded18a.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-09-03 23:05:46"

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $urlMd1, "", "", $name1)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $urlMd2, "", "", $name2)

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("I2").Value = $name1
$wsDe.Range("J2").Value = "07f5c7a9-97d9-4a11-9d2a-0d22f6f50138.bdc57f7efbf1e5bb7004563fe238e320ad0e99a9.de-de.xlf"
$wsDe.Range("K2").Value = "2016-09-03 23:05:54"

$wsDe.Range("I3").Value = $name2
$wsDe.Range("J3").Value = "7f7f1aed-0a42-45ff-a0aa-1e0ff486b802.6f682d1f368a32d8bd4f8857abd3b3d30eded18a.de-de.xlf"
$wsDe.Range("K3").Value = "2016-09-03 23:05:54"

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $urlMd1, "", "", $name1)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $urlMd2, "", "", $name2)

Write-Output "Handback report generated."
